# Auto-generated Excel COM-interop script applying the diff changes
# described in the commit "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 21
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 29
$ws.Range("AD2").Value = 81
$ws.Range("AE2").Value = 17

# Row 13
$ws.Range("J13").Value = 1.04
$ws.Range("K13").Value = 13

# Row 22
$ws.Range("G22").Value = 3.3
$ws.Range("H22").Value = 3.7
$ws.Range("I22").Value = 1.98
$ws.Range("K22").Value = 9.5
$ws.Range("M22").Value = 4.9
$ws.Range("O22").Value = 2.57
$ws.Range("P22").Value = 1.26
$ws.Range("Q22").Value = 3.5
$ws.Range("S22").Value = 2.67
$ws.Range("T22").Value = 17.5
$ws.Range("U22").Value = 25
$ws.Range("V22").Value = 11.75
$ws.Range("W22").Value = 50
$ws.Range("X22").Value = 24
$ws.Range("Y22").Value = 22
$ws.Range("Z22").Value = 9.5
$ws.Range("AA22").Value = 8
$ws.Range("AB22").Value = 10.75
$ws.Range("AC22").Value = 30
$ws.Range("AE22").Value = 12
$ws.Range("AF22").Value = 13
$ws.Range("AG22").Value = 8.75
$ws.Range("AH22").Value = 20
$ws.Range("AI22").Value = 13.5
$ws.Range("AJ22").Value = 17

# Row 28
$ws.Range("G28").Value = 1.27
$ws.Range("H28").Value = 5
$ws.Range("I28").Value = 8.25
$ws.Range("T28").Value = 7.1
$ws.Range("U28").Value = 5.9
$ws.Range("V28").Value = 7.5
$ws.Range("W28").Value = 6.8
$ws.Range("X28").Value = 8.75
$ws.Range("Y28").Value = 20
$ws.Range("Z28").Value = 15.5
$ws.Range("AA28").Value = 9
$ws.Range("AB28").Value = 17
$ws.Range("AC28").Value = 65
$ws.Range("AD28").Value = 400
$ws.Range("AE28").Value = 20
$ws.Range("AF28").Value = 50
$ws.Range("AG28").Value = 21
$ws.Range("AH28").Value = 150
$ws.Range("AI28").Value = 70
$ws.Range("AJ28").Value = 55

# Row 29
$ws.Range("G29").Value = 1.98
$ws.Range("H29").Value = 3.1
$ws.Range("I29").Value = 3.8
$ws.Range("L29").Value = 1.52
$ws.Range("M29").Value = 2.2
$ws.Range("N29").Value = 2.47
$ws.Range("O29").Value = 1.42
$ws.Range("Q29").Value = 2.12
$ws.Range("R29").Value = 2.2
$ws.Range("S29").Value = 1.52
$ws.Range("T29").Value = 5.1
$ws.Range("U29").Value = 7.8
$ws.Range("W29").Value = 17
$ws.Range("Z29").Value = 6.3
$ws.Range("AB29").Value = 22
$ws.Range("AC29").Value = 150
$ws.Range("AE29").Value = 7.8
$ws.Range("AF29").Value = 18
$ws.Range("AG29").Value = 14.5
$ws.Range("AH29").Value = 60
$ws.Range("AI29").Value = 50
$ws.Range("AJ29").Value = 75

# Row 40
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 4.1
$ws.Range("I40").Value = 1.85
$ws.Range("R40").Value = 1.25
$ws.Range("S40").Value = 3.75
$ws.Range("T40").Value = 26
$ws.Range("U40").Value = 29
$ws.Range("V40").Value = 13
$ws.Range("W40").Value = 41
$ws.Range("X40").Value = 21
$ws.Range("Y40").Value = 19
$ws.Range("AA40").Value = 11
$ws.Range("AB40").Value = 11
$ws.Range("AC40").Value = 21
$ws.Range("AE40").Value = 21
$ws.Range("AF40").Value = 19
$ws.Range("AG40").Value = 11
$ws.Range("AH40").Value = 21
$ws.Range("AI40").Value = 13

# Row 46
$ws.Range("G46").Value = 1.95
$ws.Range("H46").Value = 3.6
$ws.Range("I46").Value = 3.7
$ws.Range("AA46").Value = 7
$ws.Range("AD46").Value = 251

# Row 48
$ws.Range("G48").Value = 1.34
$ws.Range("H48").Value = 5.1
$ws.Range("I48").Value = 7.2
$ws.Range("J48").Value = 1.02
$ws.Range("K48").Value = 10
$ws.Range("M48").Value = 5.3
$ws.Range("N48").Value = 1.39
$ws.Range("O48").Value = 2.77
$ws.Range("P48").Value = 1.23
$ws.Range("Q48").Value = 3.75
$ws.Range("R48").Value = 1.62
$ws.Range("S48").Value = 2.15
$ws.Range("T48").Value = 10.75
$ws.Range("U48").Value = 8.75
$ws.Range("W48").Value = 10
$ws.Range("Y48").Value = 19
$ws.Range("Z48").Value = 10
$ws.Range("AA48").Value = 10.75
$ws.Range("AB48").Value = 16.5
$ws.Range("AE48").Value = 27

# Row 52
$ws.Range("H52").Value = 3.25
$ws.Range("I52").Value = 3.6
$ws.Range("J52").Value = 1.1
$ws.Range("K52").Value = 7
$ws.Range("N52").Value = 2.4
$ws.Range("O52").Value = 1.53
$ws.Range("P52").Value = 1.53
$ws.Range("Q52").Value = 2.38
$ws.Range("R52").Value = 2.2
$ws.Range("S52").Value = 1.62
$ws.Range("T52").Value = 5.5
$ws.Range("Y52").Value = 41
$ws.Range("Z52").Value = 7
$ws.Range("AA52").Value = 6.5
$ws.Range("AB52").Value = 21
$ws.Range("AC52").Value = 81
$ws.Range("AI52").Value = 41
